# Insert a new weekly price record for Perejil (Feria Lagunitas de Puerto Montt)
# as row 63. All existing rows from 63..157 shift down by one (to 64..158),
# which is exactly what Excel's native row Insert does (shifts cells down
# and carries the row-above formatting along, matching style "2" on column D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 63, pushing old row 63..157 down to 64..158
$ws.Rows("63:63").Insert()

# Populate the newly inserted row 63 with the new record's data
$ws.Range("A63").Value = 4
$ws.Range("B63").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C63").Value = "Los Lagos"
$ws.Range("D63").Value = 44477
$ws.Range("E63").Value = 10
$ws.Range("F63").Value = 100112044
$ws.Range("G63").Value = "Perejil"
$ws.Range("H63").Value = "Sin especificar"
$ws.Range("I63").Value = "Primera"
$ws.Range("J63").Value = 180
$ws.Range("K63").Value = 4500
$ws.Range("L63").Value = 4500
$ws.Range("M63").Value = 4500
$ws.Range("N63").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O63").Value = "Región Metropolitana"
$ws.Range("P63").Value = 1500
$ws.Range("Q63").Value = 3
$ws.Range("R63").Value = "Hortaliza"
